# Refresh market-price snapshot cells (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) for each job sheet's Leve table, as produced by the scheduled
# market-data runner. Row/column coordinates and values below mirror the runner's
# latest Universalis pull for the Phoenix data center.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 266.16666
$ws.Range("I12").Value = 266.16666
$ws.Range("K12").Value = 266.16666
$ws.Range("M12").Value = -96.16665999999998
$ws.Range("H20").Value = 3566.6667
$ws.Range("I20").Value = 3566.6667
$ws.Range("K20").Value = 3566.6667
$ws.Range("M20").Value = -3336.6667
$ws.Range("H33").Value = 2350.6
$ws.Range("I33").Value = 2498.5
$ws.Range("J33").Value = 280
$ws.Range("K33").Value = 2498.5
$ws.Range("L33").Value = 280
$ws.Range("M33").Value = -2269.5
$ws.Range("N33").Value = -738
$ws.Range("H34").Value = 5500
$ws.Range("I34").Value = 5500
$ws.Range("K34").Value = 5500
$ws.Range("M34").Value = -5297
$ws.Range("H35").Value = 3566.6667
$ws.Range("I35").Value = 3566.6667
$ws.Range("K35").Value = 3566.6667
$ws.Range("M35").Value = -3187.6667
$ws.Range("H36").Value = 5500
$ws.Range("I36").Value = 5500
$ws.Range("K36").Value = 5500
$ws.Range("M36").Value = -4785
$ws.Range("H41").Value = 1728.1818
$ws.Range("I41").Value = 2028.8889
$ws.Range("J41").Value = 375
$ws.Range("K41").Value = 2028.8889
$ws.Range("L41").Value = 375
$ws.Range("M41").Value = -1588.8889
$ws.Range("N41").Value = -1255
$ws.Range("H42").Value = 1207.125
$ws.Range("I42").Value = 1322.5714
$ws.Range("K42").Value = 3967.7142
$ws.Range("M42").Value = -3737.7142
$ws.Range("H86").Value = 1305.75
$ws.Range("I86").Value = 1628.4
$ws.Range("J86").Value = 768
$ws.Range("K86").Value = 1628.4
$ws.Range("L86").Value = 768
$ws.Range("M86").Value = -505.4000000000001
$ws.Range("N86").Value = -3014
$ws.Range("H89").Value = 1305.75
$ws.Range("I89").Value = 1628.4
$ws.Range("J89").Value = 768
$ws.Range("K89").Value = 8142
$ws.Range("L89").Value = 3840
$ws.Range("M89").Value = -2526
$ws.Range("N89").Value = -15072
$ws.Range("H94").Value = 897.8333
$ws.Range("I94").Value = 897.8333
$ws.Range("K94").Value = 897.8333
$ws.Range("M94").Value = -446.8333
$ws.Range("H96").Value = 461.95834
$ws.Range("I96").Value = 358.33334
$ws.Range("J96").Value = 565.5833
$ws.Range("K96").Value = 1075.00002
$ws.Range("L96").Value = 1696.7499
$ws.Range("M96").Value = 297.9999800000001
$ws.Range("N96").Value = -4442.7499
$ws.Range("H100").Value = 2461.6155
$ws.Range("J100").Value = 1599.1666
$ws.Range("L100").Value = 1599.1666
$ws.Range("N100").Value = -2681.1666
$ws.Range("H103").Value = 1377.7
$ws.Range("I103").Value = 721.3333
$ws.Range("J103").Value = 1493.5294
$ws.Range("K103").Value = 2163.9999
$ws.Range("L103").Value = 4480.5882
$ws.Range("M103").Value = -1577.9999
$ws.Range("N103").Value = -5652.5882
$ws.Range("H111").Value = 1149.7778
$ws.Range("J111").Value = 1298.5
$ws.Range("L111").Value = 3895.5
$ws.Range("N111").Value = -10029.5
$ws.Range("H131").Value = 2299.8
$ws.Range("J131").Value = 3033
$ws.Range("L131").Value = 9099
$ws.Range("N131").Value = -19179
$ws.Range("H132").Value = 1576.2808
$ws.Range("I132").Value = 1282.7142
$ws.Range("J132").Value = 3374.375
$ws.Range("K132").Value = 3848.1426
$ws.Range("L132").Value = 10123.125
$ws.Range("M132").Value = -1318.1426
$ws.Range("N132").Value = -15183.125
$ws.Range("H137").Value = 2097.825
$ws.Range("I137").Value = 1374.52
$ws.Range("K137").Value = 4123.559999999999
$ws.Range("M137").Value = -1573.559999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21228.584
$ws.Range("I32").Value = 19328.543
$ws.Range("J32").Value = 133331
$ws.Range("K32").Value = 19328.543
$ws.Range("L32").Value = 133331
$ws.Range("M32").Value = -19041.543
$ws.Range("N32").Value = -133905
$ws.Range("H56").Value = 20999.334
$ws.Range("I56").Value = 14999
$ws.Range("J56").Value = 23999.5
$ws.Range("K56").Value = 14999
$ws.Range("L56").Value = 23999.5
$ws.Range("M56").Value = -14257
$ws.Range("N56").Value = -25483.5
$ws.Range("H80").Value = 35067.332
$ws.Range("J80").Value = 40080.8
$ws.Range("L80").Value = 40080.8
$ws.Range("N80").Value = -42076.8
$ws.Range("H83").Value = 35067.332
$ws.Range("J83").Value = 40080.8
$ws.Range("L83").Value = 120242.4
$ws.Range("N83").Value = -130226.4
$ws.Range("H88").Value = 2020
$ws.Range("I88").Value = 1560
$ws.Range("J88").Value = 2250
$ws.Range("K88").Value = 1560
$ws.Range("L88").Value = 2250
$ws.Range("M88").Value = -1154
$ws.Range("N88").Value = -3062
$ws.Range("H91").Value = 2020
$ws.Range("I91").Value = 1560
$ws.Range("J91").Value = 2250
$ws.Range("K91").Value = 1560
$ws.Range("L91").Value = 2250
$ws.Range("M91").Value = -156
$ws.Range("N91").Value = -5058
$ws.Range("I97").Value = 951.5333000000001
$ws.Range("J97").Value = 142858020
$ws.Range("K97").Value = 951.5333000000001
$ws.Range("L97").Value = 142858020
$ws.Range("M97").Value = -455.5333000000001
$ws.Range("N97").Value = -142859012

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 302264.66
$ws.Range("I86").Value = 1799
$ws.Range("J86").Value = 452497.5
$ws.Range("K86").Value = 1799
$ws.Range("L86").Value = 452497.5
$ws.Range("M86").Value = -676
$ws.Range("N86").Value = -454743.5
$ws.Range("H89").Value = 302264.66
$ws.Range("I89").Value = 1799
$ws.Range("J89").Value = 452497.5
$ws.Range("K89").Value = 8995
$ws.Range("L89").Value = 2262487.5
$ws.Range("M89").Value = -3379
$ws.Range("N89").Value = -2273719.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 30201.5
$ws.Range("I93").Value = 16117.286
$ws.Range("J93").Value = 44285.715
$ws.Range("K93").Value = 16117.286
$ws.Range("L93").Value = 44285.715
$ws.Range("M93").Value = -14245.286
$ws.Range("N93").Value = -48029.715
$ws.Range("H99").Value = 3829.4614
$ws.Range("I99").Value = 3503.652
$ws.Range("J99").Value = 6327.3335
$ws.Range("K99").Value = 3503.652
$ws.Range("L99").Value = 6327.3335
$ws.Range("M99").Value = -2005.652
$ws.Range("N99").Value = -9323.333500000001
$ws.Range("H126").Value = 3829.4614
$ws.Range("I126").Value = 3503.652
$ws.Range("J126").Value = 6327.3335
$ws.Range("K126").Value = 10510.956
$ws.Range("L126").Value = 18982.0005
$ws.Range("M126").Value = -8040.956
$ws.Range("N126").Value = -23922.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3016.2222
$ws.Range("I129").Value = 1066
$ws.Range("J129").Value = 3991.3333
$ws.Range("K129").Value = 3198
$ws.Range("L129").Value = 11973.9999
$ws.Range("M129").Value = 1802
$ws.Range("N129").Value = -21973.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 52632190
$ws.Range("I97").Value = 62500430
$ws.Range("J97").Value = 1589.3334
$ws.Range("K97").Value = 62500430
$ws.Range("L97").Value = 1589.3334
$ws.Range("M97").Value = -62499934
$ws.Range("N97").Value = -2581.3334
$ws.Range("H132").Value = 3529.9666
$ws.Range("I132").Value = 3742.423
$ws.Range("K132").Value = 11227.269
$ws.Range("M132").Value = -8697.269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5449.091
$ws.Range("I61").Value = 3743.625
$ws.Range("J61").Value = 9997
$ws.Range("K61").Value = 3743.625
$ws.Range("L61").Value = 9997
$ws.Range("M61").Value = -3541.625
$ws.Range("N61").Value = -10401
$ws.Range("H93").Value = 1571.5264
$ws.Range("I93").Value = 1324.4286
$ws.Range("J93").Value = 2263.4
$ws.Range("K93").Value = 1324.4286
$ws.Range("L93").Value = 2263.4
$ws.Range("M93").Value = -76.42859999999996
$ws.Range("N93").Value = -4759.4
$ws.Range("H113").Value = 5449.091
$ws.Range("I113").Value = 3743.625
$ws.Range("J113").Value = 9997
$ws.Range("K113").Value = 3743.625
$ws.Range("L113").Value = 9997
$ws.Range("M113").Value = -1573.625
$ws.Range("N113").Value = -14337
$ws.Range("H132").Value = 4282.9414
$ws.Range("I132").Value = 3347.7
$ws.Range("J132").Value = 5619
$ws.Range("K132").Value = 10043.1
$ws.Range("L132").Value = 16857
$ws.Range("M132").Value = -7513.099999999999
$ws.Range("N132").Value = -21917

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1340.7894
$ws.Range("I107").Value = 768.73334
$ws.Range("K107").Value = 2306.20002
$ws.Range("M107").Value = -386.2000200000002
$ws.Range("H126").Value = 25775.125
$ws.Range("I126").Value = 31261.77
$ws.Range("K126").Value = 93785.31
$ws.Range("M126").Value = -91315.31
$ws.Range("H132").Value = 3776.5173
$ws.Range("I132").Value = 2019.2727
$ws.Range("K132").Value = 6057.8181
$ws.Range("M132").Value = -3527.8181
